$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Widen column B (39 -> 52 characters in the saved file) ---
# NOTE: the ColumnWidth property is expressed in "characters" and Excel's
# internal pixel-rounding means the value that lands on an exact file width
# of 52 is not simply 52; empirically 51.17 rounds to a stored width of 52.
$ws.Columns.Item(2).ColumnWidth = 51.17

# --- Insert a new row above row 2, shifting existing rows 2-12 down to 3-13 ---
$ws.Rows.Item(2).Insert()

# --- Remove all existing hyperlinks; row-insert does not renumber the
#     worksheet's <hyperlinks> refs automatically in this runtime, so we
#     rebuild them all from scratch at their correct (shifted) locations ---
$ws.Cells.Hyperlinks.Delete()

# --- Fill in the brand-new row 2 with the newly scraped listing ---
#     (column A's timestamp is set below, together with all the other rows)
$ws.Range("B2").Value = "最新AI活用、書き伝票から在庫更新請求入金消込までの完全自動化スキーム構築Claude/Gemini"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5490911"
$ws.Range("G2").Value = 395
$ws.Range("H2").Value = "🔥AI,Ai ◆自動化"

# --- Refresh the "retrieved at" timestamp for every row (2-13) ---
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-02-13 02:37:42"
}

# --- Re-create the hyperlinks for column F, rows 2-13, matching the Hyperlink style ---
$urls = @(
    "https://www.lancers.jp/work/detail/5490911",
    "https://www.lancers.jp/work/detail/5489981",
    "https://www.lancers.jp/work/detail/5490828",
    "https://www.lancers.jp/work/detail/5490408",
    "https://www.lancers.jp/work/detail/5490638",
    "https://www.lancers.jp/work/detail/5477084",
    "https://www.lancers.jp/work/detail/5490679",
    "https://www.lancers.jp/work/detail/5490478",
    "https://www.lancers.jp/work/detail/5486471",
    "https://www.lancers.jp/work/detail/5490062",
    "https://www.lancers.jp/work/detail/5490905",
    "https://www.lancers.jp/work/detail/5490407"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}

Write-Host "done"
